# Updates the cryptos list (Price / Volume(1h) columns, and row 51 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.164.71"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "1.627.44"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.525"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.12%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Value = "1.637.85"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.545"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.86"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "27.161.58"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -0.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("E22").Value = "  -5.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.78%  "
$ws.Range("E27").Value = "  -0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.36"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("D33").Value = "1.352.30"
$ws.Range("E33").Value = "  +6.90%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("E35").Value = "  -0.36%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.551"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.805"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.39%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.93%  "
$ws.Range("E43").Value = "  -0.49%  "
$ws.Range("D44").Value = "1.764.96"
$ws.Range("E44").Value = "  -1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  +1.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.855"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +28.68%  "
$ws.Range("E48").Value = "  +2.71%  "
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.22%  "
